$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("B3. Mô hình quan hệ")
$ws3.Range("Q7").Copy() | Out-Null
try {
  $ws3.Range("Q8").PasteSpecial(-4144) | Out-Null
  Write-Host "pasted comments"
} catch {
  Write-Host "ERR: $_"
}
$excel.CutCopyMode = $false
Write-Host "Q8 has comment:" $ws3.Range("Q8").Comment
if ($ws3.Range("Q8").Comment) {
  Write-Host "text:" $ws3.Range("Q8").Comment.Text()
}
